# Auto-generated Excel COM-interop script to apply market-price data updates
# to the Shiva_Profits leve-crafting workbook, as produced by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 25836.666
$ws.Range("I16").Value = 13750
$ws.Range("J16").Value = 50010
$ws.Range("K16").Value = 13750
$ws.Range("L16").Value = 50010
$ws.Range("M16").Value = -13520
$ws.Range("N16").Value = -50470

$ws.Range("H32").Value = 15386427
$ws.Range("J32").Value = 2155
$ws.Range("L32").Value = 2155
$ws.Range("N32").Value = -2807

$ws.Range("H62").Value = 3422.5
$ws.Range("J62").Value = 2895
$ws.Range("L62").Value = 2895
$ws.Range("N62").Value = -4143

$ws.Range("H65").Value = 3422.5
$ws.Range("J65").Value = 2895
$ws.Range("L65").Value = 14475
$ws.Range("N65").Value = -20715

$ws.Range("H76").Value = 3666.6667
$ws.Range("J76").Value = 3500
$ws.Range("L76").Value = 3500
$ws.Range("N76").Value = -4130

$ws.Range("H79").Value = 3666.6667
$ws.Range("J79").Value = 3500
$ws.Range("L79").Value = 3500
$ws.Range("N79").Value = -5684

$ws.Range("H113").Value = 10800.111
$ws.Range("I113").Value = 12733
$ws.Range("J113").Value = 6934.3335
$ws.Range("K113").Value = 12733
$ws.Range("L113").Value = 6934.3335
$ws.Range("M113").Value = -9479
$ws.Range("N113").Value = -13442.3335

$ws.Range("H116").Value = 4440.125
$ws.Range("J116").Value = 4539.8
$ws.Range("L116").Value = 4539.8
$ws.Range("N116").Value = -11423.8

$ws.Range("H138").Value = 3317.46
$ws.Range("J138").Value = 3316.8057
$ws.Range("L138").Value = 9950.417099999999
$ws.Range("N138").Value = -20230.4171

$ws.Range("H140").Value = 67991.375
$ws.Range("J140").Value = 67991.375
$ws.Range("L140").Value = 67991.375
$ws.Range("N140").Value = -78351.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1588853.5
$ws.Range("I2").Value = 1906350.1
$ws.Range("J2").Value = 1371
$ws.Range("K2").Value = 1906350.1
$ws.Range("L2").Value = 1371
$ws.Range("M2").Value = -1906237.1
$ws.Range("N2").Value = -1597

$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H88").Value = 1620.0385
$ws.Range("J88").Value = 1696.5714
$ws.Range("L88").Value = 1696.5714
$ws.Range("N88").Value = -2508.5714

$ws.Range("H91").Value = 1620.0385
$ws.Range("J91").Value = 1696.5714
$ws.Range("L91").Value = 1696.5714
$ws.Range("N91").Value = -4504.5714

$ws.Range("H102").Value = 7672.4287
$ws.Range("I102").Value = 8703.4
$ws.Range("J102").Value = 7099.6665
$ws.Range("K102").Value = 8703.4
$ws.Range("L102").Value = 7099.6665
$ws.Range("M102").Value = -7081.4
$ws.Range("N102").Value = -10343.6665

$ws.Range("H116").Value = 1588853.5
$ws.Range("I116").Value = 1906350.1
$ws.Range("J116").Value = 1371
$ws.Range("K116").Value = 1906350.1
$ws.Range("L116").Value = 1371
$ws.Range("M116").Value = -1904056.1
$ws.Range("N116").Value = -5959

$ws.Range("H122").Value = 9995
$ws.Range("I122").Value = 9994.286
$ws.Range("K122").Value = 29982.858
$ws.Range("M122").Value = -27532.858

$ws.Range("H138").Value = 45982.8
$ws.Range("J138").Value = 45982.8
$ws.Range("L138").Value = 45982.8
$ws.Range("N138").Value = -56262.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1588853.5
$ws.Range("I3").Value = 1906350.1
$ws.Range("J3").Value = 1371
$ws.Range("K3").Value = 1906350.1
$ws.Range("L3").Value = 1371
$ws.Range("M3").Value = -1906236.1
$ws.Range("N3").Value = -1599

$ws.Range("H94").Value = 2352.5
$ws.Range("I94").Value = 603
$ws.Range("K94").Value = 603
$ws.Range("M94").Value = -152

$ws.Range("H105").Value = 2810.524
$ws.Range("I105").Value = 1760.4839
$ws.Range("J105").Value = 5769.727
$ws.Range("K105").Value = 1760.4839
$ws.Range("L105").Value = 5769.727
$ws.Range("M105").Value = -13.48389999999995
$ws.Range("N105").Value = -9263.726999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 51.23077
$ws.Range("I7").Value = 52.166668
$ws.Range("K7").Value = 52.166668
$ws.Range("M7").Value = 60.833332

$ws.Range("H16").Value = 7949.3335
$ws.Range("I16").Value = 7399
$ws.Range("K16").Value = 7399
$ws.Range("M16").Value = -7112

$ws.Range("H103").Value = 22596.125
$ws.Range("I103").Value = 22596.125
$ws.Range("K103").Value = 22596.125
$ws.Range("M103").Value = -21424.125

$ws.Range("H107").Value = 3306.3572
$ws.Range("I107").Value = 2279
$ws.Range("K107").Value = 2279
$ws.Range("M107").Value = -359

$ws.Range("H113").Value = 7949.3335
$ws.Range("I113").Value = 7399
$ws.Range("K113").Value = 7399
$ws.Range("M113").Value = -5229

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("N19").ClearContents()

$ws.Range("H68").Value = 26999.75
$ws.Range("J68").Value = 26999.75
$ws.Range("L68").Value = 80999.25
$ws.Range("N68").Value = -82621.25

$ws.Range("H71").Value = 26999.75
$ws.Range("J71").Value = 26999.75
$ws.Range("L71").Value = 242997.75
$ws.Range("N71").Value = -251109.75

$ws.Range("H104").Value = 13893.591
$ws.Range("I104").Value = 10194.714
$ws.Range("J104").Value = 15619.733
$ws.Range("K104").Value = 30584.142
$ws.Range("L104").Value = 46859.199
$ws.Range("M104").Value = -27963.142
$ws.Range("N104").Value = -52101.199

$ws.Range("H114").Value = 5557198.5
$ws.Range("I114").Value = 665.125
$ws.Range("J114").Value = 11907522
$ws.Range("K114").Value = 1995.375
$ws.Range("L114").Value = 35722566
$ws.Range("M114").Value = 1258.625
$ws.Range("N114").Value = -35729074

$ws.Range("H121").Value = 8336648.5
$ws.Range("J121").Value = 3633.7222
$ws.Range("L121").Value = 10901.1666
$ws.Range("N121").Value = -13521.1666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()

$ws.Range("H122").Value = 772201.9399999999
$ws.Range("I122").Value = 2003365.2
$ws.Range("K122").Value = 6010095.6
$ws.Range("M122").Value = -6007645.6

$ws.Range("H138").Value = 249997.67
$ws.Range("J138").Value = 249997.67
$ws.Range("L138").Value = 249997.67
$ws.Range("N138").Value = -260277.67

$ws.Range("H139").Value = 379500
$ws.Range("J139").Value = 379500
$ws.Range("L139").Value = 379500
$ws.Range("N139").Value = -389780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 7981.6113
$ws.Range("I61").Value = 8284.941000000001
$ws.Range("K61").Value = 8284.941000000001
$ws.Range("M61").Value = -8082.941000000001

$ws.Range("H93").Value = 2039.7
$ws.Range("I93").Value = 2313.6875
$ws.Range("K93").Value = 2313.6875
$ws.Range("M93").Value = -1065.6875

$ws.Range("H113").Value = 7981.6113
$ws.Range("I113").Value = 8284.941000000001
$ws.Range("K113").Value = 8284.941000000001
$ws.Range("M113").Value = -6114.941000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 61916.234
$ws.Range("I81").Value = 1967.7
$ws.Range("J81").Value = 147557
$ws.Range("K81").Value = 3935.4
$ws.Range("L81").Value = 295114
$ws.Range("M81").Value = -2874.4
$ws.Range("N81").Value = -297236

$ws.Range("H84").Value = 61916.234
$ws.Range("I84").Value = 1967.7
$ws.Range("J84").Value = 147557
$ws.Range("K84").Value = 19677
$ws.Range("L84").Value = 1475570
$ws.Range("M84").Value = -14373
$ws.Range("N84").Value = -1486178

$ws.Range("H132").Value = 1797.6364
$ws.Range("I132").Value = 1259.2222
$ws.Range("K132").Value = 3777.6666
$ws.Range("M132").Value = -1247.6666
